$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old extra columns (F:M) that are no longer part of the table
$ws.Range("F2:M2").Clear()

# Give row 4 (2017, stored as text) the bold/bordered formatting up front by
# copying it from A2, then switch the whole row to Text format so values
# such as "1,061" are preserved verbatim instead of being parsed as numbers.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A4:E4").NumberFormat = "@"

# Build the sheet column by column so the shared-string table is populated
# in the same header-then-2017-value order as the source workbook.
$ws.Range("A1").Value = "Year"
$ws.Range("A4").Value = "2017"

$ws.Range("B1").Value = "Total participants"
$ws.Range("B4").Value = "1,061"

$ws.Range("C1").Value = "Total for company sports activities"
$ws.Range("C4").Value = "740"

$ws.Range("D1").Value = "Total for sickness/injury prevention courses"
$ws.Range("D4").Value = "321"

$ws.Range("E1").Value = "Average health rate for the year (%)"
$ws.Range("E4").Value = "96.5"

# Row 2 - 2019 data (numbers); A2 already carries the bold/bordered style
$ws.Range("A2").Value = 2019
$ws.Range("B2").Value = 784
$ws.Range("C2").Value = 554
$ws.Range("D2").Value = 230
$ws.Range("E2").Value = 95.3

# Row 3 - 2018 data (numbers); copy A2's direct formatting onto A3 first
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A3").Value = 2018
$ws.Range("B3").Value = 991
$ws.Range("C3").Value = 683
$ws.Range("D3").Value = 308
$ws.Range("E3").Value = 95.3

# Re-apply the bold/bordered format to A4 (the text assignment above reset its xf
# to the plain Text style), while B4:E4 fall back to the default style.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B4:E4").Style = "Normal"

$excel.CutCopyMode = 0
